$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new talk entry: IRB Board Training (row 12/13)
$ws.Range("A12").Value = 5
$ws.Range("C12").Value = "April 26, 2022"
$ws.Range("F12").Value = "Invited by"
$ws.Range("F13").Value = "Topic: Power and Sample Size Considerations in Human Subjects Research"
$ws.Range("B12").Value = "IRB Board Training"
$ws.Range("D12").Value = "Institutional Review Board"
$ws.Range("E12").Value = "Utah State University"

# Update the existing HDFS guest-lecture entry (row 10/11): shorten the course title
$ws.Range("D10").Value = "HDFS 7200: Meta Analysis"

$ws.Range("D12").Select()
